$d = $word.ActiveDocument

$d.Content.Find.Execute("25+61=", $true, $false, $false, $false, $false, $true, 1, $false, "59-27=", 2) | Out-Null
$d.Content.Find.Execute("25+49=", $true, $false, $false, $false, $false, $true, 1, $false, "23+49=", 2) | Out-Null
$d.Content.Find.Execute("91-2=", $true, $false, $false, $false, $false, $true, 1, $false, "21+54=", 2) | Out-Null
$d.Content.Find.Execute("46-18=", $true, $false, $false, $false, $false, $true, 1, $false, "93-63=", 2) | Out-Null
$d.Content.Find.Execute("62-26=", $true, $false, $false, $false, $false, $true, 1, $false, "70-43=", 2) | Out-Null
$d.Content.Find.Execute("64+34=", $true, $false, $false, $false, $false, $true, 1, $false, "70+11=", 2) | Out-Null
$d.Content.Find.Execute("49-40=", $true, $false, $false, $false, $false, $true, 1, $false, "70+6=", 2) | Out-Null
$d.Content.Find.Execute("90-60=", $true, $false, $false, $false, $false, $true, 1, $false, "42-4=", 2) | Out-Null
$d.Content.Find.Execute("14+23=", $true, $false, $false, $false, $false, $true, 1, $false, "66+15=", 2) | Out-Null
$d.Content.Find.Execute("46+5=", $true, $false, $false, $false, $false, $true, 1, $false, "62+32=", 2) | Out-Null
$d.Content.Find.Execute("11+58=", $true, $false, $false, $false, $false, $true, 1, $false, "85-65=", 2) | Out-Null
$d.Content.Find.Execute("57-52=", $true, $false, $false, $false, $false, $true, 1, $false, "37+60=", 2) | Out-Null
$d.Content.Find.Execute("51-26=", $true, $false, $false, $false, $false, $true, 1, $false, "70-25=", 2) | Out-Null
$d.Content.Find.Execute("56-13=", $true, $false, $false, $false, $false, $true, 1, $false, "40+24=", 2) | Out-Null
$d.Content.Find.Execute("26+55=", $true, $false, $false, $false, $false, $true, 1, $false, "69-49=", 2) | Out-Null
$d.Content.Find.Execute("37+55=", $true, $false, $false, $false, $false, $true, 1, $false, "63-36=", 2) | Out-Null
$d.Content.Find.Execute("64-43=", $true, $false, $false, $false, $false, $true, 1, $false, "37-28=", 2) | Out-Null
$d.Content.Find.Execute("76-12=", $true, $false, $false, $false, $false, $true, 1, $false, "80-37=", 2) | Out-Null
$d.Content.Find.Execute("62-40=", $true, $false, $false, $false, $false, $true, 1, $false, "65+2=", 2) | Out-Null
$d.Content.Find.Execute("47-23=", $true, $false, $false, $false, $false, $true, 1, $false, "4+35=", 2) | Out-Null
$d.Content.Find.Execute("20-12=", $true, $false, $false, $false, $false, $true, 1, $false, "69-59=", 2) | Out-Null
$d.Content.Find.Execute("88-36=", $true, $false, $false, $false, $false, $true, 1, $false, "70-58=", 2) | Out-Null
$d.Content.Find.Execute("58+33=", $true, $false, $false, $false, $false, $true, 1, $false, "9+1=", 2) | Out-Null
$d.Content.Find.Execute("0+78=", $true, $false, $false, $false, $false, $true, 1, $false, "44+44=", 2) | Out-Null
$d.Content.Find.Execute("61-45=", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=", 2) | Out-Null
$d.Content.Find.Execute("34+50=", $true, $false, $false, $false, $false, $true, 1, $false, "66-28=", 2) | Out-Null
$d.Content.Find.Execute("93-88=", $true, $false, $false, $false, $false, $true, 1, $false, "47-29=", 2) | Out-Null
$d.Content.Find.Execute("34-11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+59=", 2) | Out-Null
$d.Content.Find.Execute("84-20=", $true, $false, $false, $false, $false, $true, 1, $false, "57+16=", 2) | Out-Null
$d.Content.Find.Execute("28-15=", $true, $false, $false, $false, $false, $true, 1, $false, "13+54=", 2) | Out-Null
$d.Content.Find.Execute("32-24=", $true, $false, $false, $false, $false, $true, 1, $false, "55+42=", 2) | Out-Null
$d.Content.Find.Execute("16+27=", $true, $false, $false, $false, $false, $true, 1, $false, "55-53=", 2) | Out-Null
$d.Content.Find.Execute("98-77=", $true, $false, $false, $false, $false, $true, 1, $false, "7+63=", 2) | Out-Null
$d.Content.Find.Execute("11+49=", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=", 2) | Out-Null
$d.Content.Find.Execute("30+53=", $true, $false, $false, $false, $false, $true, 1, $false, "16+65=", 2) | Out-Null
$d.Content.Find.Execute("96-77=", $true, $false, $false, $false, $false, $true, 1, $false, "38-1=", 2) | Out-Null
$d.Content.Find.Execute("62-15=", $true, $false, $false, $false, $false, $true, 1, $false, "63+14=", 2) | Out-Null
$d.Content.Find.Execute("6+72=", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=", 2) | Out-Null
$d.Content.Find.Execute("31+49=", $true, $false, $false, $false, $false, $true, 1, $false, "12+15=", 2) | Out-Null
$d.Content.Find.Execute("92-81=", $true, $false, $false, $false, $false, $true, 1, $false, "43+17=", 2) | Out-Null
$d.Content.Find.Execute("33+60=", $true, $false, $false, $false, $false, $true, 1, $false, "32-2=", 2) | Out-Null
$d.Content.Find.Execute("35+27=", $true, $false, $false, $false, $false, $true, 1, $false, "39-30=", 2) | Out-Null
$d.Content.Find.Execute("96-24=", $true, $false, $false, $false, $false, $true, 1, $false, "53-0=", 2) | Out-Null
$d.Content.Find.Execute("18+76=", $true, $false, $false, $false, $false, $true, 1, $false, "11+84=", 2) | Out-Null
$d.Content.Find.Execute("22-11=", $true, $false, $false, $false, $false, $true, 1, $false, "31+18=", 2) | Out-Null
$d.Content.Find.Execute("73-9=", $true, $false, $false, $false, $false, $true, 1, $false, "32-23=", 2) | Out-Null
$d.Content.Find.Execute("34-17=", $true, $false, $false, $false, $false, $true, 1, $false, "13+56=", 2) | Out-Null
$d.Content.Find.Execute("7+25=", $true, $false, $false, $false, $false, $true, 1, $false, "37-18=", 2) | Out-Null
$d.Content.Find.Execute("54-14=", $true, $false, $false, $false, $false, $true, 1, $false, "6+39=", 2) | Out-Null
$d.Content.Find.Execute("69-64=", $true, $false, $false, $false, $false, $true, 1, $false, "2+16=", 2) | Out-Null
$d.Content.Find.Execute("28+24=", $true, $false, $false, $false, $false, $true, 1, $false, "79-57=", 2) | Out-Null
$d.Content.Find.Execute("14+12=", $true, $false, $false, $false, $false, $true, 1, $false, "34-29=", 2) | Out-Null
$d.Content.Find.Execute("80-70=", $true, $false, $false, $false, $false, $true, 1, $false, "73-52=", 2) | Out-Null
$d.Content.Find.Execute("85-36=", $true, $false, $false, $false, $false, $true, 1, $false, "21+15=", 2) | Out-Null
$d.Content.Find.Execute("36-24=", $true, $false, $false, $false, $false, $true, 1, $false, "24-14=", 2) | Out-Null
$d.Content.Find.Execute("84-15=", $true, $false, $false, $false, $false, $true, 1, $false, "17+24=", 2) | Out-Null
$d.Content.Find.Execute("45+37=", $true, $false, $false, $false, $false, $true, 1, $false, "48+24=", 2) | Out-Null
$d.Content.Find.Execute("77-25=", $true, $false, $false, $false, $false, $true, 1, $false, "44-37=", 2) | Out-Null
$d.Content.Find.Execute("81+6=", $true, $false, $false, $false, $false, $true, 1, $false, "43+12=", 2) | Out-Null
$d.Content.Find.Execute("26+65=", $true, $false, $false, $false, $false, $true, 1, $false, "67-17=", 2) | Out-Null
$d.Content.Find.Execute("11+57=", $true, $false, $false, $false, $false, $true, 1, $false, "82-25=", 2) | Out-Null
$d.Content.Find.Execute("35+2=", $true, $false, $false, $false, $false, $true, 1, $false, "39+49=", 2) | Out-Null
$d.Content.Find.Execute("14+28=", $true, $false, $false, $false, $false, $true, 1, $false, "96-34=", 2) | Out-Null
$d.Content.Find.Execute("8+58=", $true, $false, $false, $false, $false, $true, 1, $false, "2+96=", 2) | Out-Null
$d.Content.Find.Execute("23+33=", $true, $false, $false, $false, $false, $true, 1, $false, "31-3=", 2) | Out-Null
$d.Content.Find.Execute("74-27=", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=", 2) | Out-Null
$d.Content.Find.Execute("97+0=", $true, $false, $false, $false, $false, $true, 1, $false, "68+7=", 2) | Out-Null
$d.Content.Find.Execute("90-22=", $true, $false, $false, $false, $false, $true, 1, $false, "15+41=", 2) | Out-Null
$d.Content.Find.Execute("91-6=", $true, $false, $false, $false, $false, $true, 1, $false, "48-28=", 2) | Out-Null
$d.Content.Find.Execute("4+42=", $true, $false, $false, $false, $false, $true, 1, $false, "26+21=", 2) | Out-Null
$d.Content.Find.Execute("39+3=", $true, $false, $false, $false, $false, $true, 1, $false, "22+36=", 2) | Out-Null
$d.Content.Find.Execute("22-14=", $true, $false, $false, $false, $false, $true, 1, $false, "44-6=", 2) | Out-Null
$d.Content.Find.Execute("67-10=", $true, $false, $false, $false, $false, $true, 1, $false, "58-55=", 2) | Out-Null
$d.Content.Find.Execute("49+10=", $true, $false, $false, $false, $false, $true, 1, $false, "41-29=", 2) | Out-Null
$d.Content.Find.Execute("44+26=", $true, $false, $false, $false, $false, $true, 1, $false, "52+38=", 2) | Out-Null
$d.Content.Find.Execute("94-68=", $true, $false, $false, $false, $false, $true, 1, $false, "29+6=", 2) | Out-Null
$d.Content.Find.Execute("73+19=", $true, $false, $false, $false, $false, $true, 1, $false, "75-57=", 2) | Out-Null
$d.Content.Find.Execute("6+38=", $true, $false, $false, $false, $false, $true, 1, $false, "97-62=", 2) | Out-Null
$d.Content.Find.Execute("63-19=", $true, $false, $false, $false, $false, $true, 1, $false, "13+60=", 2) | Out-Null
$d.Content.Find.Execute("87-33=", $true, $false, $false, $false, $false, $true, 1, $false, "77-21=", 2) | Out-Null
$d.Content.Find.Execute("93-92=", $true, $false, $false, $false, $false, $true, 1, $false, "3+75=", 2) | Out-Null
$d.Content.Find.Execute("57-10=", $true, $false, $false, $false, $false, $true, 1, $false, "49+25=", 2) | Out-Null
$d.Content.Find.Execute("90-3=", $true, $false, $false, $false, $false, $true, 1, $false, "86-0=", 2) | Out-Null
$d.Content.Find.Execute("86-40=", $true, $false, $false, $false, $false, $true, 1, $false, "59+20=", 2) | Out-Null
$d.Content.Find.Execute("36-20=", $true, $false, $false, $false, $false, $true, 1, $false, "42+19=", 2) | Out-Null
$d.Content.Find.Execute("62-45=", $true, $false, $false, $false, $false, $true, 1, $false, "9+68=", 2) | Out-Null
$d.Content.Find.Execute("10+70=", $true, $false, $false, $false, $false, $true, 1, $false, "32-3=", 2) | Out-Null
$d.Content.Find.Execute("94-57=", $true, $false, $false, $false, $false, $true, 1, $false, "15+34=", 2) | Out-Null
$d.Content.Find.Execute("56-7=", $true, $false, $false, $false, $false, $true, 1, $false, "15+68=", 2) | Out-Null
$d.Content.Find.Execute("91-47=", $true, $false, $false, $false, $false, $true, 1, $false, "71-56=", 2) | Out-Null
$d.Content.Find.Execute("56-6=", $true, $false, $false, $false, $false, $true, 1, $false, "83+0=", 2) | Out-Null
$d.Content.Find.Execute("69-51=", $true, $false, $false, $false, $false, $true, 1, $false, "87-64=", 2) | Out-Null
$d.Content.Find.Execute("27+29=", $true, $false, $false, $false, $false, $true, 1, $false, "18+67=", 2) | Out-Null
$d.Content.Find.Execute("50-33=", $true, $false, $false, $false, $false, $true, 1, $false, "47+48=", 2) | Out-Null
$d.Content.Find.Execute("64-39=", $true, $false, $false, $false, $false, $true, 1, $false, "72+10=", 2) | Out-Null
$d.Content.Find.Execute("35+24=", $true, $false, $false, $false, $false, $true, 1, $false, "18+6=", 2) | Out-Null
$d.Content.Find.Execute("53+33=", $true, $false, $false, $false, $false, $true, 1, $false, "21+54=", 2) | Out-Null
$d.Content.Find.Execute("75-17=", $true, $false, $false, $false, $false, $true, 1, $false, "65+23=", 2) | Out-Null
$d.Content.Find.Execute("95-44=", $true, $false, $false, $false, $false, $true, 1, $false, "43+32=", 2) | Out-Null
$d.Content.Find.Execute("13-2=", $true, $false, $false, $false, $false, $true, 1, $false, "54+31=", 2) | Out-Null
